# Database Design: Update Relationship Movie-Category
# Appends ": <type>" (and a couple of Vietnamese clarifications) after the
# field-name paragraphs in the "Database:" outline, and reshapes the
# Movie.Description/Duration/Rated/PublishedYear paragraphs to also carry
# their type annotations.

$d = $word.ActiveDocument

function Append-AfterText([int]$paraIndex, [string]$oldText, [string]$suffix) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, ($oldText + $suffix), 2)
}

# User
Append-AfterText 56 "UserId" ": int"
Append-AfterText 57 "Email" ": string"
Append-AfterText 58 "Username" ": string"

# Password: append after the existing run without disturbing the
# lastRenderedPageBreak marker that precedes the "Password" text.
$p = $d.Paragraphs(59)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Collapse(0)
$r.InsertAfter(": string")

# Role (User.RoleId + the Role entity itself)
Append-AfterText 60 "RoleId" ": int"
Append-AfterText 62 "RoleId" ": int"

# RoleName paragraph also has a trailing <w:br/> run after the text -
# scope the Find to the paragraph so the break is left untouched.
Append-AfterText 63 "RoleName" ": string"

# Movie
Append-AfterText 65 "MovieId" ": int"
Append-AfterText 66 "MovieName" ": string"
Append-AfterText 67 "VideoPath" ": string (đường dẫn đến video của phim)"
Append-AfterText 68 "Description" ": string"

# Duration / Rated / PublishedYear get restructured with the type
# annotation inserted right after the field name.
$p = $d.Paragraphs(69)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Text = "Duration: int (thời lượng phim – tính theo phút)"

$p = $d.Paragraphs(70)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Text = "Rated: int (độ tuổi xem)"

$p = $d.Paragraphs(71)
$r = $p.Range
$r.MoveEnd(1, -1)
$r.Text = "PublishedYear: int"

Append-AfterText 72 "Country" ": string"
Append-AfterText 73 "ImagePath" ": string (đường dẫn đến ảnh poster của phim)"

# Category
Append-AfterText 75 "CategoryId" ": string"
Append-AfterText 76 "CategoryName" ": name"
Append-AfterText 77 "Category-Movie" " (tự động gen trong DB, ko có trong Models)"

# paragraph 78 (CategoryId under Category-Movie) is left unchanged.
